# Other procedures add in
$wb = $excel.ActiveWorkbook

# --- Array Types sheet: populate header + one data row ---
$wsArray = $wb.Worksheets.Item("Array Types")
$wsArray.Range("A1").Value = "Package Name"
$wsArray.Range("B1").Value = "Type Name"
$wsArray.Range("C1").Value = "Index Type"
$wsArray.Range("D1").Value = "Array Values Type"
$wsArray.Range("A2").Value = "General_Checks"
$wsArray.Range("B2").Value = "Mask_Type"
$wsArray.Range("C2").Value = "Mask_Index_Array_Type"
$wsArray.Range("D2").Value = "Boolean"

# Copy the header formatting (fill + border) from an existing header row
$wb.Worksheets.Item("Range Types").Range("A1:D1").Copy()
$wsArray.Range("A1:D1").PasteSpecial(-4122)

$wsArray.Range("D4").Select()

# --- Constants sheet: populate header row ---
$wsConst = $wb.Worksheets.Item("Constants")
$wsConst.Range("A1").Value = "Package Name"
$wsConst.Range("B1").Value = "Constant Name"
$wsConst.Range("C1").Value = "Type"
$wsConst.Range("E1").Value = "Value"
$wsConst.Range("D1").Value = "Attribute"

# Copy the header formatting (fill + border) from an existing header row
$wb.Worksheets.Item("Range Types").Range("A1:E1").Copy()
$wsConst.Range("A1:E1").PasteSpecial(-4122)

$wsConst.Range("D2").Select()

# --- Range Types sheet: update values + selection ---
$wsRange = $wb.Worksheets.Item("Range Types")
$wsRange.Range("F2").Value = 13
$wsRange.Range("E3").Value = 0
$wsRange.Range("F3").Value = 13
$wsRange.Range("A1:D1").Select()

# --- Procedures sheet: update selection ---
$wsProc = $wb.Worksheets.Item("Procedures")
$wsProc.Activate()
$wsProc.Range("B7").Select()
